# V3 Gerbers and STEP Files
# The BOM title cell (merged D1:F3) is renamed from the V2 part number to
# the V3 part number, and the sheet's saved selection follows the edit
# onto that merged title cell (mirrors what Excel records when the user
# clicks/edits D1, which is merged into D1:F3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "T-1D18W547848A BOM"
$ws.Range("D1:F3").Select() | Out-Null
